# Repull data, push all data, update mean calculation for the dSF column (F).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new dSF (column F) value, per repulled data.
$updates = @{
    9  = -2
    13 = 3
    20 = -3
    28 = 1
    32 = -1
    35 = -2
    36 = 0
    38 = -1
    49 = 3
    53 = 1
    56 = -3
    57 = -2
    61 = -2
    70 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
